$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 47337.332
$ws.Range("I13").Value = 2000
$ws.Range("J13").Value = 70006
$ws.Range("K13").Value = 2000
$ws.Range("L13").Value = 70006
$ws.Range("M13").Value = -1831
$ws.Range("N13").Value = -70344

$ws.Range("H40").Value = 1826
$ws.Range("I40").Value = 842.8570999999999
$ws.Range("J40").Value = 2138.818
$ws.Range("K40").Value = 842.8570999999999
$ws.Range("L40").Value = 2138.818
$ws.Range("M40").Value = -667.8570999999999
$ws.Range("N40").Value = -2488.818

$ws.Range("H92").Value = 905286.5600000001
$ws.Range("I92").Value = 961804.5
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 961804.5
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -960556.5
$ws.Range("N92").Value = -3496

$ws.Range("H106").Value = 3652
$ws.Range("I106").Value = 3726.5
$ws.Range("J106").Value = 3503
$ws.Range("K106").Value = 3726.5
$ws.Range("L106").Value = 3503
$ws.Range("M106").Value = -3095.5
$ws.Range("N106").Value = -4765

$ws.Range("H111").Value = 6478.7646
$ws.Range("I111").Value = 5064.143
$ws.Range("J111").Value = 7469
$ws.Range("K111").Value = 15192.429
$ws.Range("L111").Value = 22407
$ws.Range("M111").Value = -12125.429
$ws.Range("N111").Value = -28541

$ws.Range("H138").Value = 1844.2
$ws.Range("I138").Value = 1496.683
$ws.Range("J138").Value = 2335.5173
$ws.Range("K138").Value = 4490.049
$ws.Range("L138").Value = 7006.5519
$ws.Range("M138").Value = 649.951
$ws.Range("N138").Value = -17286.5519

$ws.Range("H141").Value = 8094.227
$ws.Range("I141").Value = 2398.4
$ws.Range("J141").Value = 65052.5
$ws.Range("K141").Value = 7195.200000000001
$ws.Range("L141").Value = 195157.5
$ws.Range("M141").Value = -2015.200000000001
$ws.Range("N141").Value = -205517.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5966750.5
$ws.Range("I32").Value = 8967.143
$ws.Range("J32").Value = 23840102
$ws.Range("K32").Value = 8967.143
$ws.Range("L32").Value = 23840102
$ws.Range("M32").Value = -8680.143
$ws.Range("N32").Value = -23840676

$ws.Range("H45").Value = 2549.818
$ws.Range("I45").Value = 1492.625
$ws.Range("J45").Value = 3544.8235
$ws.Range("K45").Value = 1492.625
$ws.Range("L45").Value = 3544.8235
$ws.Range("M45").Value = -1115.625
$ws.Range("N45").Value = -4298.8235

$ws.Range("H74").Value = 46877976
$ws.Range("I74").Value = 75003624
$ws.Range("J74").Value = 1900.5834
$ws.Range("K74").Value = 75003624
$ws.Range("L74").Value = 1900.5834
$ws.Range("M74").Value = -75002750
$ws.Range("N74").Value = -3648.5834

$ws.Range("H77").Value = 46877976
$ws.Range("I77").Value = 75003624
$ws.Range("J77").Value = 1900.5834
$ws.Range("K77").Value = 375018120
$ws.Range("L77").Value = 9502.916999999999
$ws.Range("M77").Value = -375013752
$ws.Range("N77").Value = -18238.917

$ws.Range("H132").Value = 1739182.8
$ws.Range("I132").Value = 1394.4546
$ws.Range("J132").Value = 4469993
$ws.Range("K132").Value = 4183.3638
$ws.Range("L132").Value = 13409979
$ws.Range("M132").Value = -1653.3638
$ws.Range("N132").Value = -13415039

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 583.4
$ws.Range("I11").Value = 477.8
$ws.Range("J11").Value = 689
$ws.Range("K11").Value = 477.8
$ws.Range("L11").Value = 689
$ws.Range("M11").Value = -337.8
$ws.Range("N11").Value = -969

$ws.Range("H75").Value = 8999.666999999999
$ws.Range("I75").Value = 7000
$ws.Range("J75").Value = 9999.5
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 9999.5
$ws.Range("M75").Value = -6064
$ws.Range("N75").Value = -11871.5

$ws.Range("H78").Value = 8999.666999999999
$ws.Range("I78").Value = 7000
$ws.Range("J78").Value = 9999.5
$ws.Range("K78").Value = 21000
$ws.Range("L78").Value = 29998.5
$ws.Range("M78").Value = -16320
$ws.Range("N78").Value = -39358.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2595
$ws.Range("I11").Value = 200
$ws.Range("J11").Value = 4990
$ws.Range("K11").Value = 200
$ws.Range("L11").Value = 4990
$ws.Range("M11").Value = -60
$ws.Range("N11").Value = -5270

$ws.Range("H51").Value = 28333.334
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 28333.334
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 28333.334
$ws.Range("N51").Value = -29805.334

$ws.Range("H59").Value = 32693.94
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 32693.94
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 32693.94
$ws.Range("N59").Value = -34983.94

$ws.Range("H60").Value = 11605
$ws.Range("I60").Value = 1480
$ws.Range("J60").Value = 14980
$ws.Range("K60").Value = 1480
$ws.Range("L60").Value = 14980
$ws.Range("M60").Value = -969
$ws.Range("N60").Value = -16002

$ws.Range("H61").Value = 28333.334
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 28333.334
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 28333.334
$ws.Range("N61").Value = -29029.334

$ws.Range("H68").Value = 26489.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 26489.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 26489.5
$ws.Range("N68").Value = -27987.5

$ws.Range("H71").Value = 26489.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 26489.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 79468.5
$ws.Range("N71").Value = -86956.5

$ws.Range("H74").Value = 39656.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 39656.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 39656.5
$ws.Range("N74").Value = -41404.5

$ws.Range("H77").Value = 39656.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 39656.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 118969.5
$ws.Range("N77").Value = -127705.5

$ws.Range("H122").Value = 55558290
$ws.Range("I122").Value = 66667812
$ws.Range("J122").Value = 10666.667
$ws.Range("K122").Value = 200003436
$ws.Range("L122").Value = 32000.001
$ws.Range("M122").Value = -200000986
$ws.Range("N122").Value = -36900.001

$ws.Range("H134").Value = 2350.6775
$ws.Range("I134").Value = 1309.0555
$ws.Range("J134").Value = 3792.923
$ws.Range("K134").Value = 3927.1665
$ws.Range("L134").Value = 11378.769
$ws.Range("M134").Value = -1392.1665
$ws.Range("N134").Value = -16448.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 53753.75
$ws.Range("I5").Value = 5000
$ws.Range("J5").Value = 70005
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 70005
$ws.Range("M5").Value = -4888
$ws.Range("N5").Value = -70229

$ws.Range("H70").Value = 5169.7915
$ws.Range("I70").Value = 4192.5
$ws.Range("J70").Value = 6538
$ws.Range("K70").Value = 4192.5
$ws.Range("L70").Value = 6538
$ws.Range("M70").Value = -3922.5
$ws.Range("N70").Value = -7078

$ws.Range("H73").Value = 5169.7915
$ws.Range("I73").Value = 4192.5
$ws.Range("J73").Value = 6538
$ws.Range("K73").Value = 4192.5
$ws.Range("L73").Value = 6538
$ws.Range("M73").Value = -3256.5
$ws.Range("N73").Value = -8410

$ws.Range("H80").Value = 2845
$ws.Range("I80").Value = 2845
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2845
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -1847

$ws.Range("H83").Value = 2845
$ws.Range("I83").Value = 2845
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 14225
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -9233

$ws.Range("H113").Value = 2200
$ws.Range("I113").Value = 1945.4546
$ws.Range("J113").Value = 2666.6667
$ws.Range("K113").Value = 1945.4546
$ws.Range("L113").Value = 2666.6667
$ws.Range("M113").Value = 224.5454
$ws.Range("N113").Value = -7006.6667

$ws.Range("H132").Value = 948847.9399999999
$ws.Range("I132").Value = 1544812.4
$ws.Range("J132").Value = 2316.1765
$ws.Range("K132").Value = 4634437.199999999
$ws.Range("L132").Value = 6948.529500000001
$ws.Range("M132").Value = -4631907.199999999
$ws.Range("N132").Value = -12008.5295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 17531.555
$ws.Range("I122").Value = 37468
$ws.Range("J122").Value = 7563.3335
$ws.Range("K122").Value = 112404
$ws.Range("L122").Value = 22690.0005
$ws.Range("M122").Value = -109954
$ws.Range("N122").Value = -27590.0005

$ws.Range("H136").Value = 16668381
$ws.Range("I136").Value = 26316998
$ws.Range("J136").Value = 2588.182
$ws.Range("K136").Value = 78950994
$ws.Range("L136").Value = 7764.545999999999
$ws.Range("M136").Value = -78948444
$ws.Range("N136").Value = -12864.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 621.65
$ws.Range("I107").Value = 477.0625
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 1431.1875
$ws.Range("L107").Value = 3600
$ws.Range("M107").Value = 488.8125
$ws.Range("N107").Value = -7440

$ws.Range("H122").Value = 8065.0435
$ws.Range("I122").Value = 9531.625
$ws.Range("J122").Value = 4712.857
$ws.Range("K122").Value = 28594.875
$ws.Range("L122").Value = 14138.571
$ws.Range("M122").Value = -26144.875
$ws.Range("N122").Value = -19038.571

$ws.Range("H132").Value = 2760.2058
$ws.Range("I132").Value = 2496.4614
$ws.Range("J132").Value = 2923.476
$ws.Range("K132").Value = 7489.3842
$ws.Range("L132").Value = 8770.428
$ws.Range("M132").Value = -4959.3842
$ws.Range("N132").Value = -13830.428

$ws.Range("H136").Value = 9269976
$ws.Range("I136").Value = 10881428
$ws.Range("J136").Value = 4123.75
$ws.Range("K136").Value = 32644284
$ws.Range("L136").Value = 12371.25
$ws.Range("M136").Value = -32641734
$ws.Range("N136").Value = -17471.25
